# Update hourly dataset regression output (cap_gen_year-4final) with refreshed
# Coef./Std.Err./t/P>|t|/CI/coef_pos values from the re-run regression.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.1162205337124879

$ws.Range("B3").Value = 0.1159475776526635
$ws.Range("H3").Value = 0.2321681113651514

$ws.Range("B4").Value = 0.13357075316607
$ws.Range("C4").Value = 0.01179703357937737
$ws.Range("D4").Value = 20.58699692357818
$ws.Range("E4").Value = 0.03033988233797874
$ws.Range("F4").Value = 0.1104444260845312
$ws.Range("G4").Value = 0.1566970802476083
$ws.Range("H4").Value = 0.2497912868785579

$ws.Range("B5").Value = 0.1228575463317561
$ws.Range("C5").Value = 0.008156771566143299
$ws.Range("D5").Value = 20.46737717312538
$ws.Range("E5").Value = 0.02866258988314367
$ws.Range("F5").Value = 0.1068683027896997
$ws.Range("G5").Value = 0.1388467898738125
$ws.Range("H5").Value = 0.239078080044244

$ws.Range("B6").Value = 0.08791383177406166
$ws.Range("C6").Value = 0.00840672631726566
$ws.Range("D6").Value = 5.306776804638185
$ws.Range("E6").Value = 0.004379439210882941
$ws.Range("F6").Value = 0.07143459776743939
$ws.Range("G6").Value = 0.1043930657806837
$ws.Range("H6").Value = 0.2041343654865496

$ws.Range("B7").Value = 0.1007887833569711
$ws.Range("C7").Value = 0.008340925827625694
$ws.Range("D7").Value = 5.521294417558253
$ws.Range("E7").Value = 0.0002768320724300602
$ws.Range("F7").Value = 0.08443854454145189
$ws.Range("G7").Value = 0.1171390221724903
$ws.Range("H7").Value = 0.217009317069459

$ws.Range("B8").Value = 0.08076945253612626
$ws.Range("C8").Value = 0.005160975709345611
$ws.Range("D8").Value = 5.300585579962204
$ws.Range("E8").Value = [double]"2.207715863492478e-41"
$ws.Range("F8").Value = 0.07065268261852077
$ws.Range("G8").Value = 0.09088622245373243
$ws.Range("H8").Value = 0.1969899862486142

$ws.Range("B9").Value = 0.0839124389506023
$ws.Range("C9").Value = 0.004715800976607847
$ws.Range("D9").Value = 6.032688380303757
$ws.Range("E9").Value = [double]"5.06214135483377e-44"
$ws.Range("F9").Value = 0.07466831928774365
$ws.Range("G9").Value = 0.09315655861346119
$ws.Range("H9").Value = 0.2001329726630902

$ws.Range("B10").Value = 0.07958893669482638
$ws.Range("C10").Value = 0.004286041305547478
$ws.Range("D10").Value = 6.300713455353969
$ws.Range("E10").Value = [double]"1.151425826643004e-45"
$ws.Range("F10").Value = 0.0711872501973065
$ws.Range("G10").Value = 0.08799062319234656
$ws.Range("H10").Value = 0.1958094704073143

$ws.Range("B11").Value = 0.03852213677568325
$ws.Range("H11").Value = 0.1547426704881712

$ws.Range("B12").Value = 0.05231807266108257
$ws.Range("H12").Value = 0.1685386063735705

$ws.Range("B13").Value = 0.05934863508450505
$ws.Range("C13").Value = 0.008941158084238263
$ws.Range("D13").Value = 8.20126881744941
$ws.Range("E13").Value = 0.03927449497572273
$ws.Range("F13").Value = 0.04182032212418055
$ws.Range("G13").Value = 0.07687694804482947
$ws.Range("H13").Value = 0.175569168796993

$ws.Range("B14").Value = 0.06494635472541754
$ws.Range("H14").Value = 0.1811668884379055

$ws.Range("B15").Value = 0.07194714387278678
$ws.Range("C15").Value = 0.008770652708367917
$ws.Range("D15").Value = 10.61273806437143
$ws.Range("E15").Value = 0.04307445696745881
$ws.Range("F15").Value = 0.05475311190399069
$ws.Range("G15").Value = 0.0891411758415829
$ws.Range("H15").Value = 0.1881676775852747

$ws.Range("B16").Value = 0.07455289909647499
$ws.Range("C16").Value = 0.008677549097905558
$ws.Range("D16").Value = 11.05878832989381
$ws.Range("E16").Value = 0.04017530756650255
$ws.Range("F16").Value = 0.05754097458746286
$ws.Range("G16").Value = 0.09156482360548718
$ws.Range("H16").Value = 0.1907734328089629

$ws.Range("B17").Value = 0.07528958210056472
$ws.Range("C17").Value = 0.008640939338038599
$ws.Range("D17").Value = 11.35585986782622
$ws.Range("E17").Value = 0.0249703542156049
$ws.Range("F17").Value = 0.05834907655275946
$ws.Range("G17").Value = 0.09223008764836987
$ws.Range("H17").Value = 0.1915101158130526

$ws.Range("B18").Value = -0.1162205337124879
$ws.Range("C18").Value = 0.01253467892277532
$ws.Range("D18").Value = -15.4380905202572
$ws.Range("E18").Value = 0.01882419988679694
$ws.Range("F18").Value = -0.1407934207662741
$ws.Range("G18").Value = -0.09164764665870184

$ws.Range("B19").Value = 0.07828611377103283
$ws.Range("C19").Value = 0.008575166622851921
$ws.Range("D19").Value = 11.87429842691712
$ws.Range("E19").Value = 0.02520922667591334
$ws.Range("F19").Value = 0.06147606037796181
$ws.Range("G19").Value = 0.0950961671641038
$ws.Range("H19").Value = 0.1945066474835208

$ws.Range("B20").Value = 0.08057695117083762
$ws.Range("C20").Value = 0.009289761127057723
$ws.Range("D20").Value = 11.92277877382736
$ws.Range("E20").Value = 0.02815314867482911
$ws.Range("F20").Value = 0.06236609839545573
$ws.Range("G20").Value = 0.0987878039462195
$ws.Range("H20").Value = 0.1967974848833255

$ws.Range("B21").Value = 0.08427720798608049
$ws.Range("C21").Value = 0.009337534876854294
$ws.Range("D21").Value = 12.65246232189531
$ws.Range("E21").Value = 0.02681128183523957
$ws.Range("F21").Value = 0.06597265414501301
$ws.Range("G21").Value = 0.102581761827148
$ws.Range("H21").Value = 0.2004977416985684

$ws.Range("B22").Value = 0.08303181943855802
$ws.Range("C22").Value = 0.009025992529508232
$ws.Range("D22").Value = 13.04893686287754
$ws.Range("E22").Value = 0.03469757206667382
$ws.Range("F22").Value = 0.06533782118847553
$ws.Range("G22").Value = 0.1007258176886404
$ws.Range("H22").Value = 0.1992523531510459

$ws.Range("B23").Value = 0.08828475537890729
$ws.Range("C23").Value = 0.008705141796380966
$ws.Range("D23").Value = 12.95078220764377
$ws.Range("E23").Value = 0.02260020994592784
$ws.Range("F23").Value = 0.0712199074564837
$ws.Range("G23").Value = 0.1053496033013308
$ws.Range("H23").Value = 0.2045052890913952

$ws.Range("B24").Value = 0.09012218200051153
$ws.Range("C24").Value = 0.009152598291821379
$ws.Range("D24").Value = 12.88489703750272
$ws.Range("E24").Value = 0.03248491425092325
$ws.Range("F24").Value = 0.07218006881263604
$ws.Range("G24").Value = 0.1080642951883869
$ws.Range("H24").Value = 0.2063427157129994

$ws.Range("B25").Value = 0.09920758254248599
$ws.Range("C25").Value = 0.009175826646502223
$ws.Range("D25").Value = 14.26386389952369
$ws.Range("E25").Value = 0.01469418477787258
$ws.Range("F25").Value = 0.08122010587443278
$ws.Range("G25").Value = 0.1171950592105392
$ws.Range("H25").Value = 0.2154281162549739

$ws.Range("B26").Value = 0.1024377258300269
$ws.Range("C26").Value = 0.009207981224574481
$ws.Range("D26").Value = 14.57300151360771
$ws.Range("E26").Value = 0.01366177734269737
$ws.Range("F26").Value = 0.08438714453380632
$ws.Range("G26").Value = 0.1204883071262474
$ws.Range("H26").Value = 0.2186582595425148

$ws.Range("B27").Value = 0.1038443713201875
$ws.Range("C27").Value = 0.009298227873855947
$ws.Range("D27").Value = 14.37158067603893
$ws.Range("E27").Value = 0.03270904977906969
$ws.Range("F27").Value = 0.08561703457399401
$ws.Range("G27").Value = 0.1220717080663811
$ws.Range("H27").Value = 0.2200649050326755

$ws.Range("B28").Value = 0.1045973719965627
$ws.Range("C28").Value = 0.009541683941681636
$ws.Range("D28").Value = 14.24236688609829
$ws.Range("E28").Value = 0.04827849949460525
$ws.Range("F28").Value = 0.08589279732843817
$ws.Range("G28").Value = 0.1233019466646873
$ws.Range("H28").Value = 0.2208179057090506

$ws.Range("B29").Value = 0.08277799444486302
$ws.Range("C29").Value = 0.004224814433532974
$ws.Range("D29").Value = 6.650108748883514
$ws.Range("E29").Value = [double]"6.268562128291161e-47"
$ws.Range("F29").Value = 0.07449632753429432
$ws.Range("G29").Value = 0.09105966135543204
$ws.Range("H29").Value = 0.1989985281573509
